$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# "Elimina EC anteriores y se agregan nuevos": the two mora-period records are
# re-issued, so the period labels on the two detail rows swap places.
$ws.Range("E16").Value = "2002"
$ws.Range("E17").Value = "2001"

# "se modifica base de datos": updated Valor Mora (amount owed) for both rows.
$ws.Range("G16").Value = 828116
$ws.Range("G17").Value = 828116

Write-Host "Applied EC update"
